# This workbook is a weekly price log. A new weekly observation is added as
# a new data row right after the existing header/initial rows, at row 6,
# pushing all the rows that used to be at 6..80 down to 7..81 (the sheet
# dimension therefore grows from A1:T80 to A1:T81). The new row re-uses the
# same market/product metadata (columns A-L) as its neighbours and carries
# its own price data (columns D, M-T).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 6; everything below shifts down by one.
$ws.Rows("6:6").Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C6").Value = "Los Lagos"
$ws.Range("D6").Value = 45282
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100101
$ws.Range("H6").Value = "Berries"
$ws.Range("I6").Value = 100101001
$ws.Range("J6").Value = "Arándano (blue)"
$ws.Range("K6").Value = "Sin especificar"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 300
$ws.Range("N6").Value = 4000
$ws.Range("O6").Value = 4000
$ws.Range("P6").Value = 4000
$ws.Range("Q6").Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("R6").Value = "Región del Maule"
$ws.Range("S6").Value = 2667
$ws.Range("T6").Value = 1.5
